$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.182.97'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.378.33'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '548.86'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").Value = '138.76'
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  -2.00%  '
$ws.Range("D9").Value = '2.379.02'
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = '5.33'
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '0.349'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '25.09'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").Value = '2.790.61'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").Value = '0.0000166'
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("D17").Value = '61.107.89'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '2.381.34'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = '10.93'
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").Value = '4.15'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '320.78'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = '6.71'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '64.37'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.70'
$ws.Range("E25").Value = '  -10.92%  '
$ws.Range("E26").Value = '  +3.33%  '
$ws.Range("D27").Value = '8.18'
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '510.60'
$ws.Range("E28").Value = '  -3.84%  '
$ws.Range("D29").Value = '0.0₃0893'
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("E31").Value = '  -3.94%  '
$ws.Range("D32").Value = '1.83'
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  -3.63%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '4.68'
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("E36").Value = '  +3.01%  '
$ws.Range("E37").Value = '  -1.54%  '
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").Value = '18.59'
$ws.Range("E39").Value = '  +2.47%  '
$ws.Range("D40").Value = '147.08'
$ws.Range("E40").Value = '  +4.95%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '41.26'
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("D43").Value = '151.79'
$ws.Range("E43").Value = '  +7.62%  '
$ws.Range("D44").Value = '2.14'
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.60'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").Value = '0.0522'
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.30'
$ws.Range("E47").Value = '  -4.43%  '
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").Value = '0.0908'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").Value = '0.0224'
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").Value = '11.41'
$ws.Range("E51").Value = '  +0.33%  '
